# Generate Report for Handoff
# Mark the "9a6476c5-cba3-44c1-94b4-5de314d23fca" file as ready for
# handoff across the Overview / zh-cn / de-de sheets, updating the
# corresponding status + timestamp cells.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row 3 is the 9a6476c5 file) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-37-20 06:37:07"

# --- zh-cn sheet (row 3 is the 9a6476c5 file) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-20 06:37:04"

# --- de-de sheet (row 3 is the 9a6476c5 file) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-20 06:37:07"
